$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    116.8559799194336,
    119.3869857788086,
    119.2846755981445,
    112.9613952636719,
    115.5965270996094,
    111.9981460571289,
    113.9242248535156,
    110.1974487304688,
    113.3455429077148,
    114.0765380859375,
    115.9169845581055,
    123.7375793457031,
    134.0122222900391,
    138.5408325195312,
    166.5549163818359,
    192.1971435546875,
    188.0327911376953,
    197.8594055175781,
    200.4571075439453,
    203.4959716796875,
    202.267822265625,
    202.3260803222656,
    195.8311920166016,
    195.4060821533203,
    188.4323425292969,
    191.9466705322266,
    192.1138305664062,
    189.0509796142578,
    191.6007690429688,
    191.7070159912109,
    205.5729522705078,
    220.3265991210938,
    218.5550384521484,
    258.0477905273438,
    265.3839111328125,
    265.9122009277344,
    232.7091217041016,
    222.2828826904297,
    201.6439971923828,
    190.1631011962891,
    175.2713012695312,
    160.0762634277344,
    163.5980072021484,
    150.0546875,
    149.896484375,
    142.3611602783203,
    141.5925750732422,
    141.4111938476562
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
